# Updates the cryptos list data (price + volume columns, and two pairs of
# rows that swapped rank order) to match the latest scrape.
#
# Every written cell in columns B:E on this sheet is stored as literal text
# (t="inlineStr"/shared string), even when the text looks like a number
# (e.g. "595.43") or like a European-style thousands-grouped number
# (e.g. "68.784.27"). Plain `.Value = "..."` assignment lets Excel's COM
# layer auto-coerce numeric-looking strings into real numbers, which would
# flip the cell's stored type. Forcing NumberFormat to Text ("@") before the
# assignment keeps it literal; resetting the Style back to "Normal"
# afterwards avoids leaving a stray number-format on the cell.

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple in-place value updates (price column D / change column E) ---
# Row 2 - Bitcoin
Set-TextValue $ws.Cells.Item(2, 4) "68.784.27"
Set-TextValue $ws.Cells.Item(2, 5) "  -0.44%  "

# Row 3 - Ethereum
Set-TextValue $ws.Cells.Item(3, 4) "3.481.10"
Set-TextValue $ws.Cells.Item(3, 5) "  -1.15%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Cells.Item(4, 5) "  -0.15%  "

# Row 5 - BNB
Set-TextValue $ws.Cells.Item(5, 4) "595.43"
Set-TextValue $ws.Cells.Item(5, 5) "  +2.98%  "

# Row 6 - Solana
Set-TextValue $ws.Cells.Item(6, 4) "168.18"
Set-TextValue $ws.Cells.Item(6, 5) "  -2.26%  "

# Row 7 - XRP
Set-TextValue $ws.Cells.Item(7, 5) "  -1.84%  "

# Row 8 - LidoStakedEther
Set-TextValue $ws.Cells.Item(8, 4) "3.475.52"
Set-TextValue $ws.Cells.Item(8, 5) "  -0.99%  "

# Row 9 - USDC
Set-TextValue $ws.Cells.Item(9, 5) "  -0.03%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Cells.Item(10, 4) "0.193"
Set-TextValue $ws.Cells.Item(10, 5) "  +2.19%  "

# Row 11 - Toncoin
Set-TextValue $ws.Cells.Item(11, 4) "6.81"
Set-TextValue $ws.Cells.Item(11, 5) "  +1.74%  "

# Row 12 - Cardano
Set-TextValue $ws.Cells.Item(12, 4) "0.573"
Set-TextValue $ws.Cells.Item(12, 5) "  -5.12%  "

# Row 13 - Avalanche
Set-TextValue $ws.Cells.Item(13, 4) "46.64"
Set-TextValue $ws.Cells.Item(13, 5) "  -1.31%  "

# Row 14 - ShibaInu
Set-TextValue $ws.Cells.Item(14, 5) "  +0.99%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Cells.Item(15, 4) "4.031.09"
Set-TextValue $ws.Cells.Item(15, 5) "  -1.38%  "

# Rows 16 & 17 swapped rank: Polkadot <-> BitcoinCash
Set-TextValue $ws.Cells.Item(16, 2) "BitcoinCash"
Set-TextValue $ws.Cells.Item(16, 3) "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Cells.Item(16, 4) "613.48"
Set-TextValue $ws.Cells.Item(16, 5) "  -11.30%  "

Set-TextValue $ws.Cells.Item(17, 2) "Polkadot"
Set-TextValue $ws.Cells.Item(17, 3) "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Cells.Item(17, 4) "8.30"
Set-TextValue $ws.Cells.Item(17, 5) "  -6.21%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Cells.Item(18, 4) "3.491.44"
Set-TextValue $ws.Cells.Item(18, 5) "  -0.97%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Cells.Item(19, 4) "68.780.39"
Set-TextValue $ws.Cells.Item(19, 5) "  -0.57%  "

# Row 20 - TRON
Set-TextValue $ws.Cells.Item(20, 5) "  -2.21%  "

# Row 21 - Chainlink
Set-TextValue $ws.Cells.Item(21, 4) "17.15"
Set-TextValue $ws.Cells.Item(21, 5) "  -1.90%  "

# Row 22 - Uniswap
Set-TextValue $ws.Cells.Item(22, 4) "11.11"
Set-TextValue $ws.Cells.Item(22, 5) "  -0.58%  "

# Row 23 - Polygon
Set-TextValue $ws.Cells.Item(23, 4) "0.870"
Set-TextValue $ws.Cells.Item(23, 5) "  -3.94%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue $ws.Cells.Item(24, 4) "15.80"
Set-TextValue $ws.Cells.Item(24, 5) "  -4.74%  "

# Row 25 - Litecoin
Set-TextValue $ws.Cells.Item(25, 4) "95.83"
Set-TextValue $ws.Cells.Item(25, 5) "  -1.88%  "

# Row 26 - PancakeSwap
Set-TextValue $ws.Cells.Item(26, 4) "3.79"
Set-TextValue $ws.Cells.Item(26, 5) "  -1.37%  "

# Row 27 - LEO
Set-TextValue $ws.Cells.Item(27, 4) "5.87"
Set-TextValue $ws.Cells.Item(27, 5) "  +2.36%  "

# Row 28 - Dai
Set-TextValue $ws.Cells.Item(28, 5) "  +0.10%  "

# Row 29 - ImmutableX
Set-TextValue $ws.Cells.Item(29, 5) "  -2.00%  "

# Row 30 - RenderToken
Set-TextValue $ws.Cells.Item(30, 4) "9.11"
Set-TextValue $ws.Cells.Item(30, 5) "  -3.43%  "

# Row 31 - EthereumClassic
Set-TextValue $ws.Cells.Item(31, 4) "32.97"
Set-TextValue $ws.Cells.Item(31, 5) "  -0.73%  "

# Row 32 - Filecoin
Set-TextValue $ws.Cells.Item(32, 4) "8.40"
Set-TextValue $ws.Cells.Item(32, 5) "  -5.11%  "

# Row 33 - Stacks
Set-TextValue $ws.Cells.Item(33, 4) "3.08"
Set-TextValue $ws.Cells.Item(33, 5) "  -3.00%  "

# Row 34 - Mantle
Set-TextValue $ws.Cells.Item(34, 4) "1.32"
Set-TextValue $ws.Cells.Item(34, 5) "  -2.59%  "

# Row 35 - NEARProtocol
Set-TextValue $ws.Cells.Item(35, 4) "6.79"
Set-TextValue $ws.Cells.Item(35, 5) "  -6.69%  "

# Row 36 - Bittensor
Set-TextValue $ws.Cells.Item(36, 4) "570.84"
Set-TextValue $ws.Cells.Item(36, 5) "  +0.56%  "

# Row 37 - Cosmos
Set-TextValue $ws.Cells.Item(37, 4) "10.69"
Set-TextValue $ws.Cells.Item(37, 5) "  -1.51%  "

# Row 38 - dogwifhat
Set-TextValue $ws.Cells.Item(38, 4) "3.49"
Set-TextValue $ws.Cells.Item(38, 5) "  -4.30%  "

# Row 39 - OKB
Set-TextValue $ws.Cells.Item(39, 4) "57.01"
Set-TextValue $ws.Cells.Item(39, 5) "  -0.44%  "

# Row 40 - Hedera
Set-TextValue $ws.Cells.Item(40, 5) "  -4.10%  "

# Row 41 - FirstDigitalUSD
Set-TextValue $ws.Cells.Item(41, 4) "0.998"
Set-TextValue $ws.Cells.Item(41, 5) "  -0.16%  "

# Row 42 - Kaspa
Set-TextValue $ws.Cells.Item(42, 5) "  -0.59%  "

# Row 43 - VeChain
Set-TextValue $ws.Cells.Item(43, 4) "0.0437"
Set-TextValue $ws.Cells.Item(43, 5) "  -1.00%  "

# Row 44 - Maker
Set-TextValue $ws.Cells.Item(44, 4) "3.389.49"
Set-TextValue $ws.Cells.Item(44, 5) "  -1.64%  "

# Row 45 - TheGraph
Set-TextValue $ws.Cells.Item(45, 5) "  -4.80%  "

# Rows 46 & 47 swapped rank: InjectiveProtocol <-> PEPE
Set-TextValue $ws.Cells.Item(46, 2) "PEPE"
Set-TextValue $ws.Cells.Item(46, 3) "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Cells.Item(46, 4) "0.0₃0695"
Set-TextValue $ws.Cells.Item(46, 5) "  -1.22%  "

Set-TextValue $ws.Cells.Item(47, 2) "InjectiveProtocol"
Set-TextValue $ws.Cells.Item(47, 3) "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Cells.Item(47, 4) "32.46"
Set-TextValue $ws.Cells.Item(47, 5) "  -2.39%  "

# Rows 48 & 49 swapped rank: ThetaToken <-> Fetch.AI
Set-TextValue $ws.Cells.Item(48, 2) "Fetch.AI"
Set-TextValue $ws.Cells.Item(48, 3) "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Cells.Item(48, 4) "2.55"
Set-TextValue $ws.Cells.Item(48, 5) "  -1.33%  "

Set-TextValue $ws.Cells.Item(49, 2) "ThetaToken"
Set-TextValue $ws.Cells.Item(49, 3) "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Cells.Item(49, 4) "2.81"
Set-TextValue $ws.Cells.Item(49, 5) "  -2.47%  "

# Row 50 - Stellar
Set-TextValue $ws.Cells.Item(50, 5) "  -4.10%  "

# Row 51 - Monero
Set-TextValue $ws.Cells.Item(51, 4) "132.29"
Set-TextValue $ws.Cells.Item(51, 5) "  -1.47%  "
